$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 0.4095395551410945
$data[0,1] = 0.1722077596581713
$data[0,2] = 0.07827828017953209
$data[0,3] = 0.1274714773545842
$data[0,4] = 0
$data[0,5] = 0.7598749259888038
$data[0,6] = 0.8568042189395513
$data[0,7] = 0.887759677257641
$data[0,8] = 0
$data[0,9] = 0.2711847320485106
$data[0,10] = 0.2092779248668464
$data[0,11] = 0.1331230801806988
$data[0,12] = 0
$data[0,13] = 3.243799351104769
$data[1,0] = 0.3731603414961455
$data[1,1] = 0.1707302140052676
$data[1,2] = 0.07104946101929954
$data[1,3] = 0.1281025825239599
$data[1,4] = 0
$data[1,5] = 0.7639047782140906
$data[1,6] = 0.8625346178171753
$data[1,7] = 0.8952137378075555
$data[1,8] = 0
$data[1,9] = 0.2369687000987426
$data[1,10] = 0.2066352816822672
$data[1,11] = 0.1261100579064092
$data[1,12] = 0
$data[1,13] = 3.264073511488633
$data[2,0] = 0.3508771670434214
$data[2,1] = 0.1698207491983226
$data[2,2] = 0.06664359688569732
$data[2,3] = 0.1285392103504286
$data[2,4] = 0
$data[2,5] = 0.7668249784200043
$data[2,6] = 0.8663906553822827
$data[2,7] = 0.9001997381931446
$data[2,8] = 0
$data[2,9] = 0.21590291334158
$data[2,10] = 0.2051131320339152
$data[2,11] = 0.1218492617434386
$data[2,12] = 0
$data[2,13] = 3.278164427096812
$data[3,0] = 0.3418107588284727
$data[3,1] = 0.1694496006505517
$data[3,2] = 0.06485642092464161
$data[3,3] = 0.1287295132376212
$data[3,4] = 0
$data[3,5] = 0.7681270929297739
$data[3,6] = 0.8680469835673037
$data[3,7] = 0.9023345069459729
$data[3,8] = 0
$data[3,9] = 0.2073046157738503
$data[3,10] = 0.2045181744161582
$data[3,11] = 0.1201244547879732
$data[3,12] = 0
$data[3,13] = 3.284319765003389
$data[4,0] = 0.34030616298881
$data[4,1] = 0.1693879403480238
$data[4,2] = 0.06456016104387174
$data[4,3] = 0.1287618608094991
$data[4,4] = 0
$data[4,5] = 0.7683500790393225
$data[4,6] = 0.8683271494837257
$data[4,7] = 0.9026952018203112
$data[4,8] = 0
$data[4,9] = 0.2058760556284369
$data[4,10] = 0.2044209142253806
$data[4,11] = 0.1198387504085439
$data[4,12] = 0
$data[4,13] = 3.285366813531866
$data[5,0] = 0.350754836089294
$data[5,1] = 0.1698157458812517
$data[5,2] = 0.06661946095576354
$data[5,3] = 0.1285417267149178
$data[5,4] = 0
$data[5,5] = 0.7668420852722022
$data[5,6] = 0.8664126490889075
$data[5,7] = 0.9002281116031483
$data[5,8] = 0
$data[5,9] = 0.2157870089589551
$data[5,10] = 0.2051050055667361
$data[5,11] = 0.1218259536506672
$data[5,12] = 0
$data[5,13] = 3.278245767000357
$data[6,0] = 0.3969852033262953
$data[6,1] = 0.1716987858098591
$data[6,2] = 0.07577902468162279
$data[6,3] = 0.1276788982144534
$data[6,4] = 0
$data[6,5] = 0.7611718782373629
$data[6,6] = 0.8587100511437313
$data[6,7] = 0.8902449380516764
$data[6,8] = 0
$data[6,9] = 0.2593991744578545
$data[6,10] = 0.2083459370159062
$data[6,11] = 0.1306956781161155
$data[6,12] = 0
$data[6,13] = 3.25044907199738
$data[7,0] = 0.4880464189762108
$data[7,1] = 0.1753724434715238
$data[7,2] = 0.09399936132956555
$data[7,3] = 0.126375841437774
$data[7,4] = 0
$data[7,5] = 0.7535910460426294
$data[7,6] = 0.8462800869599221
$data[7,7] = 0.8739128155628819
$data[7,8] = 0
$data[7,9] = 0.3444515786183615
$data[7,10] = 0.2154957004453877
$data[7,11] = 0.1484433241219065
$data[7,12] = 0
$data[7,13] = 3.208966934348837
$data[8,0] = 0.5551704813136666
$data[8,1] = 0.1780586183940613
$data[8,2] = 0.1075438698172633
$data[8,3] = 0.1256545075811033
$data[8,4] = 0
$data[8,5] = 0.7501804104248748
$data[8,6] = 0.8387738831236504
$data[8,7] = 0.8638894558869268
$data[8,8] = 0
$data[8,9] = 0.4066334494632429
$data[8,10] = 0.2212300811618917
$data[8,11] = 0.1616934788158204
$data[8,12] = 0
$data[8,13] = 3.186426560104735
$data[9,0] = 0.5857502170521514
$data[9,1] = 0.1792775808566418
$data[9,2] = 0.1137401013791219
$data[9,3] = 0.1253773771131854
$data[9,4] = 0
$data[9,5] = 0.749098087601098
$data[9,6] = 0.8357112669879285
$data[9,7] = 0.8597581402813823
$data[9,8] = 0
$data[9,9] = 0.4348514523795188
$data[9,10] = 0.2239428001373369
$data[9,11] = 0.1677661641006836
$data[9,12] = 0
$data[9,13] = 3.177894794435758
$data[10,0] = 0.5973358357184679
$data[10,1] = 0.1797387130357748
$data[10,2] = 0.1160914320472131
$data[10,3] = 0.1252797506037027
$data[10,4] = 0
$data[10,5] = 0.7487557333558783
$data[10,6] = 0.8346020750451544
$data[10,7] = 0.8582552823505623
$data[10,8] = 0
$data[10,9] = 0.4455265332057081
$data[10,10] = 0.224984945669533
$data[10,11] = 0.17007210852843
$data[10,12] = 0
$data[10,13] = 3.174911552998253
$data[11,0] = 0.5948404212266212
$data[11,1] = 0.1796394210487762
$data[11,2] = 0.1155848114252933
$data[11,3] = 0.125300451091567
$data[11,4] = 0
$data[11,5] = 0.7488264630337937
$data[11,6] = 0.8348387119426235
$data[11,7] = 0.8585762110942667
$data[11,8] = 0
$data[11,9] = 0.4432279379416002
$data[11,10] = 0.2247598395932613
$data[11,11] = 0.1695752019086143
$data[11,12] = 0
$data[11,13] = 3.175543038452048
$data[12,0] = 0.5867032618313033
$data[12,1] = 0.1793155279309175
$data[12,2] = 0.1139334477913394
$data[12,3] = 0.1253691987859131
$data[12,4] = 0
$data[12,5] = 0.7490685689915608
$data[12,6] = 0.8356190002691051
$data[12,7] = 0.8596332651601699
$data[12,8] = 0
$data[12,9] = 0.4357299106596884
$data[12,10] = 0.2240282399336451
$data[12,11] = 0.1679557490594377
$data[12,12] = 0
$data[12,13] = 3.177644400402471
$data[13,0] = 0.5817197434568016
$data[13,1] = 0.1791170728675979
$data[13,2] = 0.1129225833973209
$data[13,3] = 0.1254122610580453
$data[13,4] = 0
$data[13,5] = 0.7492256569374973
$data[13,6] = 0.8361035312303784
$data[13,7] = 0.8602887605623586
$data[13,8] = 0
$data[13,9] = 0.4311357729170311
$data[13,10] = 0.2235820519220653
$data[13,11] = 0.1669646103197238
$data[13,12] = 0
$data[13,13] = 3.178963782088232
$data[14,0] = 0.5531728637210733
$data[14,1] = 0.1779788933661806
$data[14,2] = 0.1071396270979221
$data[14,3] = 0.1256736435720516
$data[14,4] = 0
$data[14,5] = 0.7502605865840337
$data[14,6] = 0.8389811100025071
$data[14,7] = 0.8641680654159387
$data[14,8] = 0
$data[14,9] = 0.4047879029989474
$data[14,10] = 0.2210548873886751
$data[14,11] = 0.1612975107181285
$data[14,12] = 0
$data[14,13] = 3.187018775994488
$data[15,0] = 0.5356712283196998
$data[15,1] = 0.1772798678080534
$data[15,2] = 0.103600842046518
$data[15,3] = 0.1258470447565667
$data[15,4] = 0
$data[15,5] = 0.7510156744773724
$data[15,6] = 0.8408365212428066
$data[15,7] = 0.866657595338971
$data[15,8] = 0
$data[15,9] = 0.3886063036477196
$data[15,10] = 0.2195311684815806
$data[15,11] = 0.1578323873748815
$data[15,12] = 0
$data[15,13] = 3.19240124551041
$data[16,0] = 0.5256089953631715
$data[16,1] = 0.1768775270888483
$data[16,2] = 0.1015687017694518
$data[16,3] = 0.1259515827205711
$data[16,4] = 0
$data[16,5] = 0.7514941420162984
$data[16,6] = 0.8419368397054114
$data[16,7] = 0.8681298274366078
$data[16,8] = 0
$data[16,9] = 0.3792926296450787
$data[16,10] = 0.2186645704519066
$data[16,11] = 0.1558435955782258
$data[16,12] = 0
$data[16,13] = 3.195659182839506
$data[17,0] = 0.5222028461593027
$data[17,1] = 0.1767412543792233
$data[17,2] = 0.1008812179430407
$data[17,3] = 0.1259878028108634
$data[17,4] = 0
$data[17,5] = 0.7516637267050896
$data[17,6] = 0.8423150814910088
$data[17,7] = 0.8686352254768224
$data[17,8] = 0
$data[17,9] = 0.3761380897841775
$data[17,10] = 0.2183728415015196
$data[17,11] = 0.1551709599526419
$data[17,12] = 0
$data[17,13] = 3.196790104298913
$data[18,0] = 0.5375338729360806
$data[18,1] = 0.1773543094538894
$data[18,2] = 0.1039772131161527
$data[18,3] = 0.1258280890270669
$data[18,4] = 0
$data[18,5] = 0.7509307233787439
$data[18,6] = 0.8406355806497032
$data[18,7] = 0.8663884077420221
$data[18,8] = 0
$data[18,9] = 0.3903295337612462
$data[18,10] = 0.2196923567427831
$data[18,11] = 0.1582008162447153
$data[18,12] = 0
$data[18,13] = 3.1918114978875
$data[19,0] = 0.5890931925700329
$data[19,1] = 0.1794106759974738
$data[19,2] = 0.1144183592975736
$data[19,3] = 0.1253488074904041
$data[19,4] = 0
$data[19,5] = 0.7489956245053264
$data[19,6] = 0.8353884391363238
$data[19,7] = 0.8593211115790815
$data[19,8] = 0
$data[19,9] = 0.4379325515270409
$data[19,10] = 0.2242427247636982
$data[19,11] = 0.1684312503744749
$data[19,12] = 0
$data[19,13] = 3.177020461451974
$data[20,0] = 0.6228233072489502
$data[20,1] = 0.1807519200935843
$data[20,2] = 0.1212710948776987
$data[20,3] = 0.1250782087223499
$data[20,4] = 0
$data[20,5] = 0.7481243677122933
$data[20,6] = 0.8322537695284495
$data[20,7] = 0.8550611622991831
$data[20,8] = 0
$data[20,9] = 0.4689826192223165
$data[20,10] = 0.2273034538911531
$data[20,11] = 0.1751543807435567
$data[20,12] = 0
$data[20,13] = 3.168796562447056
$data[21,0] = 0.6048181062131448
$data[21,1] = 0.1800363314586946
$data[21,2] = 0.1176110376192696
$data[21,3] = 0.1252187368448929
$data[21,4] = 0
$data[21,5] = 0.7485533646450335
$data[21,6] = 0.8338998624502949
$data[21,7] = 0.8573019399711548
$data[21,8] = 0
$data[21,9] = 0.4524164050546915
$data[21,10] = 0.225661967614883
$data[21,11] = 0.1715627842918721
$data[21,12] = 0
$data[21,13] = 3.173053808075082
$data[22,0] = 0.5366917722185178
$data[22,1] = 0.1773206558221574
$data[22,2] = 0.1038070484190428
$data[22,3] = 0.1258366438123844
$data[22,4] = 0
$data[22,5] = 0.7509689916018516
$data[22,6] = 0.8407263211693987
$data[22,7] = 0.8665099798400178
$data[22,8] = 0
$data[22,9] = 0.3895504945802486
$data[22,10] = 0.2196194542191847
$data[22,11] = 0.1580342390959188
$data[22,12] = 0
$data[22,13] = 3.192077613540533
$data[23,0] = 0.4633712940616022
$data[23,1] = 0.1743807878101507
$data[23,2] = 0.08904255652797133
$data[23,3] = 0.1266868261313405
$data[23,4] = 0
$data[23,5] = 0.7552628519774487
$data[23,6] = 0.8493568359845796
$data[23,7] = 0.8779839415336141
$data[23,8] = 0
$data[23,9] = 0.3214950393530671
$data[23,10] = 0.2134767509219486
$data[23,11] = 0.1436047129444049
$data[23,12] = 0
$data[23,13] = 3.21879483704501

$ws.Range("B2:O25").Value2 = $data
